$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add formulas dividing B by A for rows 4 and 5 into column C
$ws.Range("C4").Formula = "=B4/A4"
$ws.Range("C5").Formula = "=B5/A5"

# Update the current selection to K14
$ws.Range("K14").Select()
